# Update countries & provincias Spain
# - Rusia, Polonia and Filipinas got refreshed (higher) case counts and moved
#   up in the (descending, by "Casos totales") ranking, displacing the
#   country that used to occupy that rank down by one row.
# - A few other rows (Rumania, Georgia, Camboya) simply got refreshed counts
#   with no change in rank.
# - The "last updated" timestamp footer was bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rusia overtakes Belgica (and pushes Brasil down too) ---------------
# Row 13 was Belgica, row 14 was Brasil, row 15 was Rusia.
# New order: row13 = Rusia (new data), row14 = Belgica (old row13 data),
# row15 = Brasil (old row14 data).
$ws.Range("A14").Value = "Belgica"
$ws.Range("B14").Value = 37183
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 8348
$ws.Range("E14").Value = 23382
$ws.Range("F14").Value = 1119
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 5453

$ws.Range("A15").Value = "Brasil"
$ws.Range("B15").Value = 36925
$ws.Range("C15").Value = 203
$ws.Range("D15").Value = 14026
$ws.Range("E15").Value = 20527
$ws.Range("F15").Value = 6634
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 2372

$ws.Range("A13").Value = "Rusia"
$ws.Range("B13").Value = 42853
$ws.Range("C13").Value = 6060
$ws.Range("D13").Value = 3291
$ws.Range("E13").Value = 39201
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 48
$ws.Range("H13").Value = 361

# --- Polonia overtakes Ecuador -------------------------------------------
# Row 29 was Ecuador, row 30 was Polonia.
# New order: row29 = Polonia (new data), row30 = Ecuador (old row29 data).
$ws.Range("A30").Value = "Ecuador"
$ws.Range("B30").Value = 9022
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 1008
$ws.Range("E30").Value = 7558
$ws.Range("F30").Value = 168
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 456

$ws.Range("A29").Value = "Polonia"
$ws.Range("B29").Value = 9082
$ws.Range("C29").Value = 340
$ws.Range("D29").Value = 1040
$ws.Range("E29").Value = 7692
$ws.Range("F29").Value = 160
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 350

# Rumania (row 31) keeps its rank, just refreshed counts.
$ws.Range("E31").Value = 6254
$ws.Range("G31").Value = 13
$ws.Range("H31").Value = 434

# --- Filipinas overtakes Indonesia ---------------------------------------
# Row 41 was Indonesia, row 42 was Filipinas.
# New order: row41 = Filipinas (new data), row42 = Indonesia (old row41 data).
$ws.Range("A42").Value = "Indonesia"
$ws.Range("B42").Value = 6248
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 631
$ws.Range("E42").Value = 5082
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 535

$ws.Range("A41").Value = "Filipinas"
$ws.Range("B41").Value = 6259
$ws.Range("C41").Value = 172
$ws.Range("D41").Value = 572
$ws.Range("E41").Value = 5278
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 12
$ws.Range("H41").Value = 409

# --- Georgia (row 110) refreshed counts, no rank change ------------------
$ws.Range("B110").Value = 394
$ws.Range("C110").Value = 6
$ws.Range("E110").Value = 304

# --- Camboya (row 135) refreshed counts, no rank change -------------------
$ws.Range("D135").Value = 105
$ws.Range("E135").Value = 17

# --- Footer timestamp ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 10:22"
